$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.431.01'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '3.944.17'
$ws.Range('E3').Value = '  +0.29%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '502.96'
$ws.Range('E5').Value = '  +3.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.95'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.626'
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.735'
$ws.Range('E9').Value = '  -0.35%  '
$ws.Range('E10').Value = '  +3.60%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000351'
$ws.Range('E11').Value = '  -1.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.72'
$ws.Range('E12').Value = '  +1.41%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.51'
$ws.Range('E13').Value = '  -1.72%  '
$ws.Range('D14').Value = '4.574.78'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').Value = '3.938.44'
$ws.Range('E15').Value = '  +0.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.26'
$ws.Range('E16').Value = '  -2.49%  '
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('E18').Value = '  +5.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '20.04'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').Value = '69.413.84'
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '437.31'
$ws.Range('E21').Value = '  -1.52%  '
$ws.Range('E22').Value = '  -1.66%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.77'
$ws.Range('E23').Value = '  -2.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '88.98'
$ws.Range('E24').Value = '  +0.40%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.06'
$ws.Range('E25').Value = '  +5.02%  '
$ws.Range('E26').Value = '  +6.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.22'
$ws.Range('E27').Value = '  -2.54%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.23'
$ws.Range('E28').Value = '  -4.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.66'
$ws.Range('E29').Value = '  -3.28%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '710.33'
$ws.Range('E30').Value = '  -1.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.48'
$ws.Range('E31').Value = '  -2.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.129'
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.89'
$ws.Range('E33').Value = '  -1.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '64.58'
$ws.Range('E34').Value = '  +5.50%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.454'
$ws.Range('E35').Value = '  +13.16%  '
$ws.Range('D36').Value = '0.0₃0892'
$ws.Range('E36').Value = '  -2.76%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.07'
$ws.Range('E37').Value = '  -2.76%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '41.09'
$ws.Range('E38').Value = '  -3.58%  '
$ws.Range('E39').Value = '  +1.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  +2.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.90'
$ws.Range('E43').Value = '  -5.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.09'
$ws.Range('E44').Value = '  -5.95%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.06'
$ws.Range('E45').Value = '  +4.11%  '
$ws.Range('E46').Value = '  +1.00%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.39'
$ws.Range('E47').Value = '  +4.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.03'
$ws.Range('E48').Value = '  +5.85%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.40'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('D50').Value = '0.0₆0349'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.10'
$ws.Range('E51').Value = '  -2.52%  '
